$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update the time_taken values on the existing "data" sheet
$dataSheet.Range("F2").Value = "2021-10-05 14:21:06.962961"
$dataSheet.Range("F3").Value = "2021-10-05 14:21:06.962969"
$dataSheet.Range("F4").Value = "2021-10-05 14:21:06.962973"
$dataSheet.Range("F5").Value = "2021-10-05 14:21:06.962975"
$dataSheet.Range("F6").Value = "2021-10-05 14:21:06.962979"

# Add the new "metadata" sheet right after "data"
$ws = $wb.Worksheets.Add($null, $dataSheet)
$ws.Name = "metadata"

# Header row
$ws.Range("B1").Value = "data_name"
$ws.Range("C1").Value = "data_id"
$ws.Range("D1").Value = "data_version"
$ws.Range("E1").Value = "data_version_created"
$ws.Range("F1").Value = "panel_query_time"
$ws.Range("G1").Value = "panel_get_request"

# Reuse the same header style as the "data" sheet's header row
$dataSheet.Range("B1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)

# Data row
$ws.Range("A2").Value = 0
$dataSheet.Range("A2").Copy()
$ws.Range("A2").PasteSpecial(-4122)

$ws.Range("B2").Value = "Inherited MMR deficiency (Lynch syndrome)"
$ws.Range("C2").Value = 503
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1.9"
$ws.Range("E2").Value = "2021-03-03T11:37:07.405035Z"
$ws.Range("F2").Value = "2021-10-05 14:21:06.959761"
$ws.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/503/?format=json"

$dataSheet.Activate()
